$d = $word.ActiveDocument

function ReplaceInParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $found = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Host "WARNING: not found in paragraph $paraIndex : $findText"
    }
    return $r
}

# After merging runs that originally carried no formatting overrides
# (an empty <w:rPr/>), the merged run loses its (empty) <w:rPr/>
# element entirely. Toggling Bold on/off on the merged run forces the
# engine to re-emit an explicit (empty) <w:rPr/>, matching the source.
function RestoreEmptyRPr($rangeToFix) {
    $rangeToFix.Bold = 1
    $rangeToFix.Bold = 0
}

$rdquo = [char]8221

# 1. "This " + "Labtainer" + " exercise..." -> single run (keeps empty <w:rPr/>)
$r5 = ReplaceInParagraph 5 "This Labtainer exercise explores the use of the denyhosts utility on a SSH server to limit SSH login attempts from an IP address." "This Labtainer exercise explores the use of the denyhosts utility on a SSH server to limit SSH login attempts from an IP address."
RestoreEmptyRPr $r5

# 2. "T" + "he lab is started..." -> single run (keeps empty <w:rPr/>)
$r9 = ReplaceInParagraph 9 "The lab is started from the labtainer working directory on your linux host, e.g., a Linux VM. From there issue the command:" "The lab is started from the labtainer working directory on your linux host, e.g., a Linux VM. From there issue the command:"
RestoreEmptyRPr $r9

# 3. "Key " + "file" + " #1" -> single run (keeps bold rPr)
ReplaceInParagraph 22 "Key file #1" "Key file #1"

# 4. "a" + "uth.log" -> single run (keeps rFonts rPr)
ReplaceInParagraph 24 "auth.log" "auth.log"

# 5. "<rdquo> with the password " + "hank21" -> single run (keeps empty <w:rPr/>)
$target31 = "$rdquo with the password hank21"
$r31 = ReplaceInParagraph 31 $target31 $target31
RestoreEmptyRPr $r31

# 6. "Key " + "file " + "#2" -> single run (keeps bold rPr)
ReplaceInParagraph 41 "Key file #2" "Key file #2"

# 7. "d" + "enyhosts.conf" -> single run (keeps rFonts rPr)
ReplaceInParagraph 43 "denyhosts.conf" "denyhosts.conf"

# 8. "No" + "te in particular the description and values for " -> single run (keeps empty <w:rPr/>)
$r50 = ReplaceInParagraph 50 "Note in particular the description and values for " "Note in particular the description and values for "
RestoreEmptyRPr $r50

# 9. "Key " + "file " + "#3" -> single run (keeps bold rPr; trailing space run stays separate)
ReplaceInParagraph 53 "Key file #3" "Key file #3"

# 10. Font change Tlwg Typo -> Tlwg Typist for the "./bot.py hank" run only
#     (scope to just the run's text, not the paragraph mark, so pPr/rPr untouched)
$p68 = $d.Paragraphs(68)
$r68 = $p68.Range
$r68run = $d.Range($r68.Start, $r68.End - 1)
$r68run.Font.Name = "Tlwg Typist"

# 11. Shape resize (wp:extent / a:ext / Shape size)
$shape = $d.Shapes(1)
$shape.Width = 504.9
$shape.Height = 36.9

# 12. styles.xml Normal style: overflowPunct false -> true
$normalStyle = $d.Styles.Item(1)
$normalStyle.ParagraphFormat.HangingPunctuation = $true

Write-Host "edits applied"
